# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.206.03'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '1.884.29'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.77'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.687'
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.56'
$ws.Range('E8').Value = '  +2.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.352'
$ws.Range('E9').Value = '  +2.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.00'
$ws.Range('E10').Value = '  +7.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0738'
$ws.Range('E11').Value = '  +1.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.72'
$ws.Range('E13').Value = '  +7.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.775'
$ws.Range('E14').Value = '  +9.60%  '
$ws.Range('D15').Value = '2.158.91'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.94'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').Value = '1.880.78'
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').Value = '35.223.03'
$ws.Range('E18').Value = '  +1.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.07'
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('D20').Value = '0.0₃0818'
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '243.27'
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.73'
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('E23').Value = '  +5.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.63'
$ws.Range('E24').Value = '  +6.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.13'
$ws.Range('E26').Value = '  -2.00%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.19'
$ws.Range('E27').Value = '  +1.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.46'
$ws.Range('E28').Value = '  +2.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.22'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.126'
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0591'
$ws.Range('E31').Value = '  +3.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.26'
$ws.Range('E32').Value = '  +2.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.86'
$ws.Range('E33').Value = '  +20.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.14'
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  -13.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.844'
$ws.Range('E37').Value = '  +3.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.92'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0715'
$ws.Range('E39').Value = '  +7.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0219'
$ws.Range('E40').Value = '  +4.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.31'
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.03'
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('D44').Value = '1.321.90'
$ws.Range('E44').Value = '  +3.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.18'
$ws.Range('E45').Value = '  +12.12%  '
$ws.Range('E46').Value = '  +2.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0807'
$ws.Range('E47').Value = '  -1.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.40'
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.23'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').Value = '2.059.11'
$ws.Range('E51').Value = '  +0.29%  '
